# Horarios Linea 141 - actualizacion 06:35:22 (scrape "131")
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# freshly scraped arrival times, merging the new rows into the existing,
# B-column-sorted table and refreshing the "Ultima actualizacion" / "Total
# filas" info cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:35:22"
$ws1.Range("A3").Value = "Total filas: 30"

$sheet1Rows = @(
    @("05:57:13","06:01","16_SANTA ANA",4,"LP1912"),
    @("05:57:13","06:09","10_OLMOS",12,"LP1912"),
    @("05:57:13","06:16","215A_EL PATO",19,"LP1912"),
    @("06:17:28","06:17","215A_EL PATO",0,"LP1912"),
    @("05:57:13","06:30","23_HERNANDEZ",33,"LP1912"),
    @("06:17:28","06:32","23_HERNANDEZ",15,"LP1912"),
    @("05:57:13","06:34","11_ETCHEVERRY",37,"LP1912"),
    @("06:35:22","06:35","11_ETCHEVERRY",0,"LP1912"),
    @("05:57:13","06:39","17X38_ROMERO",42,"LP1912"),
    @("05:57:13","06:41","16_SANTA ANA",44,"LP1912"),
    @("05:57:13","06:57","215A_EL PATO",60,"LP1912"),
    @("05:57:13","06:59","225_GOMEZ",62,"LP1912"),
    @("06:17:28","07:15","215C_EL PATO",58,"LP1912"),
    @("05:57:13","07:16","215C_EL PATO",79,"LP1912"),
    @("05:57:13","07:19","14_ABASTO",82,"LP1912"),
    @("05:57:13","07:21","23_HERNANDEZ",84,"LP1912"),
    @("06:17:28","07:21","16_SANTA ANA",64,"LP1912"),
    @("05:57:13","07:29","17X38_ROMERO",92,"LP1912"),
    @("05:57:13","07:35","10_OLMOS",98,"LP1912"),
    @("06:17:28","07:36","27_EL RETIRO",79,"LP1912"),
    @("05:57:13","07:37","27_EL RETIRO",100,"LP1912"),
    @("06:35:22","07:44","215A_EL PATO",69,"LP1912"),
    @("05:57:13","07:55","14_ABASTO",118,"LP1912"),
    @("06:17:28","08:00","17_ROMERO",103,"LP1912"),
    @("06:17:28","08:01","16_SANTA ANA",104,"LP1912"),
    @("06:35:22","08:06","23_HERNANDEZ",91,"LP1912"),
    @("06:17:28","08:11","10_OLMOS",114,"LP1912"),
    @("06:17:28","08:13","15X38_ABASTO",116,"LP1912"),
    @("06:35:22","08:29","11_ETCHEVERRY",114,"LP1912"),
    @("06:35:22","08:29","15_ABASTO",114,"LP1912")
)

$r = 6
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:35:22"
$ws2.Range("A3").Value = "Total filas: 6"

$sheet2Rows = @(
    @("05:57:13","06:16","215A_EL PATO",19,"LP1912"),
    @("06:17:28","06:17","215A_EL PATO",0,"LP1912"),
    @("05:57:13","06:57","215A_EL PATO",60,"LP1912"),
    @("06:17:28","07:15","215C_EL PATO",58,"LP1912"),
    @("05:57:13","07:16","215C_EL PATO",79,"LP1912"),
    @("06:35:22","07:44","215A_EL PATO",69,"LP1912")
)

$r = 6
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:35:22"
